$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the ranges we will update so numeric-looking
# strings (prices, percentages, hour values) are preserved verbatim, just
# like the original inline-string cells, instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "287.22"
$ws.Range("E2").Value = "1.80%"
$ws.Range("G2").Value = "10"

$ws.Range("E3").Value = "4.17%"
$ws.Range("G3").Value = "10"

$ws.Range("D4").Value = "5.104"
$ws.Range("E4").Value = "1.61%"
$ws.Range("G4").Value = "10"

$ws.Range("D5").Value = "0.06714"
$ws.Range("E5").Value = "3.40%"
$ws.Range("G5").Value = "10"

$ws.Range("D6").Value = "7.341"
$ws.Range("E6").Value = "1.79%"
$ws.Range("G6").Value = "10"

$ws.Range("D7").Value = "3.407"
$ws.Range("E7").Value = "1.37%"
$ws.Range("G7").Value = "10"

$ws.Range("E8").Value = "-0.57%"
$ws.Range("G8").Value = "10"

$ws.Range("D9").Value = "0.9183"
$ws.Range("E9").Value = "0.02%"
$ws.Range("G9").Value = "10"

$ws.Range("D10").Value = "0.1602"
$ws.Range("E10").Value = "4.59%"
$ws.Range("G10").Value = "10"

$ws.Range("D11").Value = "0.06779"
$ws.Range("E11").Value = "6.32%"
$ws.Range("G11").Value = "10"

$ws.Range("D12").Value = "0.07778"
$ws.Range("E12").Value = "2.59%"
$ws.Range("G12").Value = "10"

$ws.Range("D13").Value = "0.02934"
$ws.Range("E13").Value = "3.46%"
$ws.Range("G13").Value = "10"

$ws.Range("D14").Value = "0.08976"
$ws.Range("E14").Value = "0.10%"
$ws.Range("G14").Value = "10"

$ws.Range("D15").Value = "0.001577"
$ws.Range("E15").Value = "-0.37%"
$ws.Range("G15").Value = "10"

$ws.Range("D16").Value = "0.04500"
$ws.Range("E16").Value = "1.53%"
$ws.Range("G16").Value = "10"

$ws.Range("D17").Value = "0.0006472"
$ws.Range("E17").Value = "1.48%"
$ws.Range("G17").Value = "10"

$ws.Range("D18").Value = "0.006219"
$ws.Range("E18").Value = "1.50%"
$ws.Range("G18").Value = "10"

$ws.Range("D19").Value = "3.447"
$ws.Range("E19").Value = "-0.03%"
$ws.Range("G19").Value = "10"

$ws.Range("D20").Value = "2.229"
$ws.Range("E20").Value = "-0.56%"
$ws.Range("G20").Value = "10"

$ws.Range("D21").Value = "0.3217"
$ws.Range("E21").Value = "1.11%"
$ws.Range("G21").Value = "10"

$ws.Range("E22").Value = "-2.22%"
$ws.Range("G22").Value = "10"

$ws.Range("D23").Value = "4.102"
$ws.Range("E23").Value = "3.39%"
$ws.Range("G23").Value = "10"

$ws.Range("E24").Value = "2.39%"
$ws.Range("G24").Value = "10"

$ws.Range("E25").Value = "0.88%"
$ws.Range("G25").Value = "10"

$ws.Range("D26").Value = "0.004122"
$ws.Range("E26").Value = "-7.45%"
$ws.Range("G26").Value = "10"

$ws.Range("E27").Value = "-0.17%"
$ws.Range("G27").Value = "10"

$ws.Range("E28").Value = "-0.22%"
$ws.Range("G28").Value = "10"

$ws.Range("G29").Value = "10"

$ws.Range("G30").Value = "10"

$ws.Range("G31").Value = "10"

$ws.Range("G32").Value = "10"

$ws.Range("G33").Value = "10"

$ws.Range("G34").Value = "10"

$ws.Range("G35").Value = "10"

$ws.Range("G36").Value = "10"

$ws.Range("G37").Value = "10"

$ws.Range("G38").Value = "10"

$ws.Range("G39").Value = "10"

$ws.Range("D40").Value = "0.04300"
$ws.Range("E40").Value = "4.80%"
$ws.Range("G40").Value = "10"

$ws.Range("D41").Value = "0.006779"
$ws.Range("E41").Value = "1.44%"
$ws.Range("G41").Value = "10"

$ws.Range("D42").Value = "0.1238"
$ws.Range("E42").Value = "0.65%"
$ws.Range("G42").Value = "10"

$ws.Range("D43").Value = "0.002226"
$ws.Range("E43").Value = "4.03%"
$ws.Range("G43").Value = "10"

$ws.Range("D44").Value = "0.01209"
$ws.Range("E44").Value = "4.68%"
$ws.Range("G44").Value = "10"

$ws.Range("D45").Value = "0.00005707"
$ws.Range("E45").Value = "5.77%"
$ws.Range("G45").Value = "10"

$ws.Range("E46").Value = "-1.29%"
$ws.Range("G46").Value = "10"

$ws.Range("E47").Value = "-29.47%"
$ws.Range("G47").Value = "10"

$ws.Range("G48").Value = "10"

$ws.Range("G49").Value = "10"

$ws.Range("G50").Value = "10"

$ws.Range("G51").Value = "10"

Write-Output "Applied cryptos price/volume/hour update for Fri Jan 13 10:10:26 UTC 2023"
